## Zeitplan.xlsx edit
## - Insert a new row (18) for "Definitive Abgabe der Maturaarbeit" with date 19.12.2011
## - Move selection to the new data-entry point (B18)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Insert a new row above the old row 18 ("?" / "Schlusspräsentation"),
# shifting it (and everything below) down by one.
$ws.Rows.Item(18).Insert()

# New row 18: submission deadline date + description
$ws.Range("A18").Value2 = 40896   # 19 December 2011 (Excel serial date)
$ws.Range("B18").Value2 = "Definitive Abgabe der Maturaarbeit"

# Update the active selection to reflect where editing left off
[void]$ws.Range("B18").Select()
